$d = $word.ActiveDocument
$wdParagraph = 4

# --- Change 1: merge the 2015-12-15 paragraph's three runs into a single run ---
$rng1 = $d.Content
$rng1.Find.Execute("2015-12-15", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Expand($wdParagraph) | Out-Null
$rng1.MoveEnd(1, -1) | Out-Null
$rng1.Delete() | Out-Null
$rng1.InsertAfter("2015-12-15 tips for running the model headless") | Out-Null

# --- Change 2: add two new dated log entries after "adding the model name..." paragraph ---
$rngAdding = $d.Content
$rngAdding.Find.Execute("adding the model name only works when you are in the address of the model",
                         $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngAdding.Expand($wdParagraph) | Out-Null

# new paragraph: 2015-12-16 ...
$rngAdding.InsertParagraphAfter() | Out-Null
$rng1216 = $rngAdding.Next($wdParagraph, 1)
$rng1216.MoveEnd(1, -1) | Out-Null
$rng1216.Text = "2015-12-16 adding the protected run in the model and start to add the algorithm for optimizing the carbon allocation"

# new paragraph: 2015-12-23 ...
$rng1216.InsertParagraphAfter() | Out-Null
$rng1223 = $rng1216.Next($wdParagraph, 1)
$rng1223.MoveEnd(1, -1) | Out-Null
$rng1223.Text = "2015-12-23 updated the carbon allocation parameter. 1. Berry mass flow set to zero when the phloem concentration is low; 2. Berry sugar inhibitor set to relates with sugar concentration. 3. Berry fraction of soluble sugar set as a function of sugar concentration in fresh weight (as we do not have good data on g/L); 4. Berry number reduce after 7 days; 5. Stop the nitrogen dynamics when optimize carbon allocation as the leaf N concentration increase too much"

# --- Change 3: add a new, truly empty paragraph right after the bookmark paragraph ---
$rngBookmark = $d.Bookmarks("_GoBack").Range
$rngBookmark.Expand($wdParagraph) | Out-Null
$rngBookmark.InsertParagraphAfter() | Out-Null
$rngBlank = $rngBookmark.Next($wdParagraph, 1)
$rngBlank.InsertAfter("x") | Out-Null
$rngBlank2 = $rngBlank
$rngBlank2.MoveEnd(1, -1) | Out-Null
$rngBlank2.Text = ""
